$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November 2022")

# Row 4 - week of 11/11/2022 (A4 already has the date, fill in the rest)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = "Sprint 5"

# Row 5 - week of 11/23/2022
$ws.Range("A5").Value = "11/23/2022"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Sprint 6"

# Row 6 - week of 11/28/2022
$ws.Range("A6").Value = "11/28/2022"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "Sprint 7"

# Row 12 - new UAT raised defect entry
$ws.Range("A12").Value = "BD-14619"
$ws.Range("B12").Value = "Medium"
$ws.Range("C12").Value = "Medium"

# Update the active selection to reflect where the user last worked
$ws.Activate()
$ws.Range("F5").Select()
